$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark several rows as passing ("x" in the Pass column, E).
$ws.Range("E5").Value = "x"
$ws.Range("E6").Value = "x"
$ws.Range("E10").Value = "x"
$ws.Range("E11").Value = "x"
$ws.Range("E13").Value = "x"
$ws.Range("E14").Value = "x"
$ws.Range("E15").Value = "x"

# Update the "Shows all X tags" text (now clarified as being "in Note") and
# capitalize "Notes" in the word-count description.
$ws.Range("C13").Value = "Shows all @ tags in Note"
$ws.Range("C14").Value = "Shows all # tags in Note"
$ws.Range("C15").Value = "Shows all ^ tags in Note"
$ws.Range("C16").Value = "Displays number of words in Notes"

# Remove the old row 18 ("15) Open URL tag in browser / Opens URL tag in browser")
# without shifting the unrelated rows further down the sheet (e.g. rows 23-27).
$ws.Range("A18:C18").ClearContents()

# Update the active selection to reflect the editor's final cursor position.
$ws.Range("D9").Select()
